$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cell H1 = "Save", matching the style of the neighboring G1 header
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the "Save" column values for each data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
